{"js": "// Replace the date heading and each \"A\u00d7B=\" multiplication prompt in the\n// practice-sheet table with the new values from the target revision.\n// Every old value is unique in the document, so a plain search+replace\n// (matchCase, no wildcards) for each pair is unambiguous and keeps the\n// original run formatting (font/size) untouched.\nconst replacements = [\n  [\"2024-08-26 Monday\", \"2024-08-27 Tuesday\"],\n  [\"383\u00d74=\", \"500\u00d76=\"],\n  [\"560\u00d73=\", \"393\u00d79=\"],\n  [\"341\u00d75=\", \"626\u00d78=\"],\n  [\"396\u00d74=\", \"368\u00d74=\"],\n  [\"205\u00d73=\", \"766\u00d72=\"],\n  [\"801\u00d79=\", \"820\u00d77=\"],\n  [\"301\u00d75=\", \"139\u00d72=\"],\n  [\"208\u00d73=\", \"586\u00d76=\"],\n  [\"999\u00d79=\", \"876\u00d79=\"],\n  [\"651\u00d78=\", \"635\u00d75=\"],\n  [\"937\u00d77=\", \"873\u00d77=\"],\n  [\"113\u00d77=\", \"455\u00d74=\"],\n  [\"622\u00d73=\", \"504\u00d75=\"],\n  [\"336\u00d79=\", \"150\u00d78=\"],\n  [\"558\u00d79=\", \"309\u00d77=\"],\n  [\"601\u00d76=\", \"933\u00d79=\"],\n  [\"574\u00d79=\", \"721\u00d75=\"],\n  [\"356\u00d78=\", \"761\u00d78=\"],\n  [\"808\u00d72=\", \"611\u00d78=\"],\n  [\"803\u00d75=\", \"721\u00d77=\"],\n  [\"885\u00d72=\", \"223\u00d75=\"],\n  [\"856\u00d78=\", \"928\u00d77=\"],\n  [\"948\u00d78=\", \"365\u00d72=\"],\n  [\"164\u00d73=\", \"448\u00d75=\"],\n  [\"648\u00d75=\", \"354\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date heading and each \"A\u00d7B=\" multiplication prompt in the\n# practice-sheet table with the new values from the target revision.\n# Every old value is unique in the document, so a plain Find/Replace\n# (MatchCase on, no wildcards) for each pair is unambiguous and keeps the\n# original run formatting (font/size) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-26 Monday\", \"2024-08-27 Tuesday\"),\n    @(\"383\u00d74=\", \"500\u00d76=\"),\n    @(\"560\u00d73=\", \"393\u00d79=\"),\n    @(\"341\u00d75=\", \"626\u00d78=\"),\n    @(\"396\u00d74=\", \"368\u00d74=\"),\n    @(\"205\u00d73=\", \"766\u00d72=\"),\n    @(\"801\u00d79=\", \"820\u00d77=\"),\n    @(\"301\u00d75=\", \"139\u00d72=\"),\n    @(\"208\u00d73=\", \"586\u00d76=\"),\n    @(\"999\u00d79=\", \"876\u00d79=\"),\n    @(\"651\u00d78=\", \"635\u00d75=\"),\n    @(\"937\u00d77=\", \"873\u00d77=\"),\n    @(\"113\u00d77=\", \"455\u00d74=\"),\n    @(\"622\u00d73=\", \"504\u00d75=\"),\n    @(\"336\u00d79=\", \"150\u00d78=\"),\n    @(\"558\u00d79=\", \"309\u00d77=\"),\n    @(\"601\u00d76=\", \"933\u00d79=\"),\n    @(\"574\u00d79=\", \"721\u00d75=\"),\n    @(\"356\u00d78=\", \"761\u00d78=\"),\n    @(\"808\u00d72=\", \"611\u00d78=\"),\n    @(\"803\u00d75=\", \"721\u00d77=\"),\n    @(\"885\u00d72=\", \"223\u00d75=\"),\n    @(\"856\u00d78=\", \"928\u00d77=\"),\n    @(\"948\u00d78=\", \"365\u00d72=\"),\n    @(\"164\u00d73=\", \"448\u00d75=\"),\n    @(\"648\u00d75=\", \"354\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $oldText\n    $range.Find.Replacement.Text = $newText\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n    $range.Find.Format = $false\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n"}
